$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Förändrad" (C) column: every existing data row's "last changed" date
#    moves from 2023-10-03 (45202) to 2023-10-04 (45203).
$ws.Range("C2:C373").Value = 45203

# 2) New cutting notification row (374) appended to the bottom of the table.
$ws.Cells.Item(374, 1).Value = "A 47082-2023"

$ws.Cells.Item(374, 2).Value = 45201
$ws.Cells.Item(374, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(374, 3).Value = 45203
$ws.Cells.Item(374, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(374, 4).Value = "NORRBOTTENS LÄN"
$ws.Cells.Item(374, 5).Value = "ÖVERTORNEÅ"

$ws.Cells.Item(374, 7).Value = 9.8
$ws.Cells.Item(374, 8).Value = 0
$ws.Cells.Item(374, 9).Value = 0
$ws.Cells.Item(374, 10).Value = 0
$ws.Cells.Item(374, 11).Value = 0
$ws.Cells.Item(374, 12).Value = 0
$ws.Cells.Item(374, 13).Value = 0
$ws.Cells.Item(374, 14).Value = 0
$ws.Cells.Item(374, 15).Value = 0
$ws.Cells.Item(374, 16).Value = 0
$ws.Cells.Item(374, 17).Value = 0

$ws.Cells.Item(374, 18).Value = ""
$ws.Cells.Item(374, 18).WrapText = $true

$ws.Rows.Item(373).RowHeight = 15
